$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Owner Name / Service Id / Program Id / Program Name block had been
# duplicated twice more across G:J and K:N. Remove those stray duplicate
# columns so the table only spans A:F again.
$ws.Range("G1:N10").EntireColumn.Delete()

# Row 10 was left without its Owner Name / Service Id / Program Id /
# Program Name values - fill them in to match the rest of the table.
$ws.Range("C7").Copy()
$ws.Range("C10").PasteSpecial(-4163)

$ws.Range("D7").Copy()
$ws.Range("D10").PasteSpecial(-4163)

$ws.Range("H1").Formula = "=""172"""
$ws.Range("H1").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("H1").ClearContents()

$ws.Range("F6").Copy()
$ws.Range("F10").PasteSpecial(-4163)

# Move the active selection to the next empty row below the table.
$ws.Range("A11").Select()
